$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$ws.Range("E2").Value = "URL"
$ws.Range("E3").Value = "https://login.salesforce.com/?locale=ca"

$ws.Range("C3").Select()
